# Daily Report update: 2026-01-21
# Adds the 2026-01-20 (serial date 46042) depository rows to Daily_Data,
# and refreshes the Today_Summary / Monthly_Stats roll-ups that change
# as a result (HSBC BANK, USA and LOOMIS INTERNATIONAL (US) LLC).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Daily_Data: append rows 244-265 for date 46042 (2026-01-20)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46042, "ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @(46042, "ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(46042, "BRINK'S, INC. Registered", 91733.761, 0, 0, 0, 0, 91733.761),
    @(46042, "BRINK'S, INC. Eligible", 27494.288, 0, 0, 0, 0, 27494.288),
    @(46042, "CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @(46042, "CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(46042, "DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @(46042, "DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @(46042, "HSBC BANK, USA Registered", 1295.223, 0, 0, 0, 99.535, 1394.758),
    @(46042, "HSBC BANK, USA Eligible", 9381.513999999999, 0, 0, 0, -99.535, 9281.978999999999),
    @(46042, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @(46042, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @(46042, "JP MORGAN CHASE BANK NA Registered", 114985.579, 0, 0, 0, 0, 114985.579),
    @(46042, "JP MORGAN CHASE BANK NA Eligible", 135413.823, 0, 0, 0, 0, 135413.823),
    @(46042, "LOOMIS INTERNATIONAL (US) LLC Registered", 52372.648, 11373.343, 0, 11373.343, 0, 63745.991),
    @(46042, "LOOMIS INTERNATIONAL (US) LLC Eligible", 132077.206, 0, 0, 0, 0, 132077.206),
    @(46042, "MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @(46042, "MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(46042, "MANFRA, TORDELLA & BROOKES, LLC Registered", 50220.42, 0, 0, 0, 0, 50220.42),
    @(46042, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 11149.237, 0, 0, 0, 0, 11149.237),
    @(46042, "STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @(46042, "STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$startRow = 244
$ws1.Range("A$startRow`:A265").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$r = $startRow
foreach ($row in $newRows) {
    $ws1.Cells.Item($r, 1).Value2 = $row[0]
    $ws1.Cells.Item($r, 2).Value2 = $row[1]
    $ws1.Cells.Item($r, 3).Value2 = $row[2]
    $ws1.Cells.Item($r, 4).Value2 = $row[3]
    $ws1.Cells.Item($r, 5).Value2 = $row[4]
    $ws1.Cells.Item($r, 6).Value2 = $row[5]
    $ws1.Cells.Item($r, 7).Value2 = $row[6]
    $ws1.Cells.Item($r, 8).Value2 = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Today_Summary: refresh rows affected by the new day's adjustments
#    (HSBC BANK, USA row 6; LOOMIS INTERNATIONAL (US) LLC row 9)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Today_Summary")

$ws2.Cells.Item(6, 2).Value2 = 9281.978999999999   # HSBC Eligible
$ws2.Cells.Item(6, 3).Value2 = 1394.758            # HSBC Registered

$ws2.Cells.Item(9, 3).Value2 = 63745.991           # LOOMIS Registered
$ws2.Cells.Item(9, 4).Value2 = 195823.197          # LOOMIS Total_Stock

# ---------------------------------------------------------------------
# 3. Monthly_Stats: refresh the month grand totals and the per-depository
#    monthly breakdown rows impacted by the new day
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Monthly_Stats")

$ws3.Cells.Item(2, 2).Value2 = 333892.192          # Eligible grand total
$ws3.Cells.Item(2, 3).Value2 = 341873.868          # Registered grand total
$ws3.Cells.Item(2, 4).Value2 = 675766.0600000001   # Grand_Total

$ws3.Cells.Item(15, 5).Value2 = 9281.978999999999  # HSBC BANK, USA Eligible - TOTAL_TODAY
$ws3.Cells.Item(16, 5).Value2 = 1394.758           # HSBC BANK, USA Registered - TOTAL_TODAY

$ws3.Cells.Item(22, 3).Value2 = 11373.343          # LOOMIS ... Registered - RECEIVED
$ws3.Cells.Item(22, 5).Value2 = 63745.991          # LOOMIS ... Registered - TOTAL_TODAY

Write-Output "Applied daily report update for 2026-01-20 (22 new rows + roll-up refresh)."
